$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (SEXO) - swap values for User3 (F) and User4 (G)
$ws.Range("F5").Value = "Hombre"
$ws.Range("G5").Value = "Mujer"

# Row 6 (EDAD)
$ws.Range("C6").Value = 23
$ws.Range("E6").Value = 57
$ws.Range("F6").Value = 26
$ws.Range("G6").Value = 24

# Row 7 (OCUPACION)
$ws.Range("C7").Value = "Estudiante"
$ws.Range("E7").Value = "Ama de casa"
$ws.Range("F7").Value = "Físico"
$ws.Range("G7").Value = "Estudiante"

# Row 8 (EXPERIENCIA TIC)
$ws.Range("C8").Value = "Alta"
$ws.Range("E8").Value = "Baja"
$ws.Range("F8").Value = "Alta"
$ws.Range("G8").Value = "Baja"

# Row 9 (PERFIL (describir))
$ws.Range("C9").Value = "Estudiante tranquilo"
$ws.Range("E9").Value = "Ama de casa feliz"
$ws.Range("F9").Value = "Trabajador tranquilo"
$ws.Range("G9").Value = "Estudiante enfadada"

# Rows 14-23 : questionnaire answers (columns C, E, F, G); D has formulas
$ws.Range("C14").Value = 4
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 3

$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2

$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 4

$ws.Range("C17").Value = 1
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 1

$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 4

$ws.Range("C19").Value = 2
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 3

$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 5
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 4

$ws.Range("C21").Value = 2
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 4

$ws.Range("C22").Value = 4
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4

$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 2

# Update view state: scroll back to top-left and move selection to G9
$ws.Range("G9").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
